{"js": "// Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table.\n// The table has 20 rows x 5 columns; every 4th row (0, 4, 8, 12, 16) holds\n// the 5 answers for that block, the rows in between are blank spacer rows.\n// Replace the text of each of those 25 cells with its new value, addressed\n// by (row, col) position so that cells sharing old/new text across\n// different positions are never confused with each other.\n\nconst table = context.document.body.tables.getFirst();\n\n// Row index (within the table) -> [old, new] text pairs for each column.\nconst updates = [\n  { row: 0, values: [\"48\u00f75=9, 3\", \"37\u00f78=4, 5\", \"10\u00f74=2, 2\", \"62\u00f72=31, 0\", \"30\u00f79=3, 3\"] },\n  { row: 4, values: [\"10\u00f73=3, 1\", \"98\u00f77=14, 0\", \"92\u00f78=11, 4\", \"72\u00f76=12, 0\", \"36\u00f79=4, 0\"] },\n  { row: 8, values: [\"10\u00f77=1, 3\", \"85\u00f73=28, 1\", \"10\u00f72=5, 0\", \"97\u00f79=10, 7\", \"51\u00f77=7, 2\"] },\n  { row: 12, values: [\"35\u00f72=17, 1\", \"36\u00f74=9, 0\", \"23\u00f74=5, 3\", \"44\u00f72=22, 0\", \"24\u00f73=8, 0\"] },\n  { row: 16, values: [\"46\u00f76=7, 4\", \"76\u00f75=15, 1\", \"92\u00f75=18, 2\", \"33\u00f79=3, 6\", \"26\u00f74=6, 2\"] },\n];\n\nfor (const { row, values } of updates) {\n  for (let col = 0; col < values.length; col++) {\n    table.getCell(row, col).value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 \"two-digit \u00f7 one-digit\" answer cells in the single table.\n# The table has 20 rows x 5 columns; every 4th row (1, 5, 9, 13, 17 in\n# Word's 1-based Cell() addressing) holds the 5 answers for that block,\n# with blank spacer rows in between. Address each cell directly by\n# (row, column) so cells that share old/new text at different positions\n# are never confused with one another.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"48\u00f75=9, 3\"\n$t.Cell(1, 2).Range.Text = \"37\u00f78=4, 5\"\n$t.Cell(1, 3).Range.Text = \"10\u00f74=2, 2\"\n$t.Cell(1, 4).Range.Text = \"62\u00f72=31, 0\"\n$t.Cell(1, 5).Range.Text = \"30\u00f79=3, 3\"\n\n$t.Cell(5, 1).Range.Text = \"10\u00f73=3, 1\"\n$t.Cell(5, 2).Range.Text = \"98\u00f77=14, 0\"\n$t.Cell(5, 3).Range.Text = \"92\u00f78=11, 4\"\n$t.Cell(5, 4).Range.Text = \"72\u00f76=12, 0\"\n$t.Cell(5, 5).Range.Text = \"36\u00f79=4, 0\"\n\n$t.Cell(9, 1).Range.Text = \"10\u00f77=1, 3\"\n$t.Cell(9, 2).Range.Text = \"85\u00f73=28, 1\"\n$t.Cell(9, 3).Range.Text = \"10\u00f72=5, 0\"\n$t.Cell(9, 4).Range.Text = \"97\u00f79=10, 7\"\n$t.Cell(9, 5).Range.Text = \"51\u00f77=7, 2\"\n\n$t.Cell(13, 1).Range.Text = \"35\u00f72=17, 1\"\n$t.Cell(13, 2).Range.Text = \"36\u00f74=9, 0\"\n$t.Cell(13, 3).Range.Text = \"23\u00f74=5, 3\"\n$t.Cell(13, 4).Range.Text = \"44\u00f72=22, 0\"\n$t.Cell(13, 5).Range.Text = \"24\u00f73=8, 0\"\n\n$t.Cell(17, 1).Range.Text = \"46\u00f76=7, 4\"\n$t.Cell(17, 2).Range.Text = \"76\u00f75=15, 1\"\n$t.Cell(17, 3).Range.Text = \"92\u00f75=18, 2\"\n$t.Cell(17, 4).Range.Text = \"33\u00f79=3, 6\"\n$t.Cell(17, 5).Range.Text = \"26\u00f74=6, 2\"\n"}
